# Applies the diff: extends the forecast table by one quarter (column BB)
# and one additional row (row 83), following the existing "Naive
# forecaster" propagation pattern.
#
# Summary of change:
#  - New column BB (col 54):
#      BB1       = 45986 (next quarter date, same style as BA1 -> s="1")
#      BB2..BB72 = copy of the corresponding BA value in that row
#      BB73..BB83 = new forecast constant 0.8783323788356512
#  - New row 83:
#      A83  = 46934 (next date in the quarterly series, style s="1")
#      BB83 = 0.8783323788356512 (see above)
#  - dimension / row "spans" attributes are refreshed automatically by
#    the application when cell values are written.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newQuarterDate = 45986
$newRowDate = 46934
$newForecastValue = 0.8783323788356512

# --- 1. Copy formatting from column BA (rows 1-82) into new column BB ---
# (This also pre-creates the BB1:BB82 cells, each inheriting BA's style,
#  i.e. BB1 gets the bordered/centered/date style "s=1", BB2:BB82 stay
#  with the default style, matching the source column.)
$ws.Range("BA1:BA82").Copy()
$ws.Range("BB1:BB82").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- 2. Fill column BB values ---
# Header cell (row 1): next quarter-end serial date.
$ws.Cells.Item(1, 54).Value = $newQuarterDate

# Rows 2-72: duplicate the existing BA value for that row.
for ($r = 2; $r -le 72; $r++) {
    $baValue = $ws.Cells.Item($r, 53).Value2
    $ws.Cells.Item($r, 54).Value = $baValue
}

# Rows 73-82: new forecast constant (diverges from the BA column value).
for ($r = 73; $r -le 82; $r++) {
    $ws.Cells.Item($r, 54).Value = $newForecastValue
}

# --- 3. Add new row 83 ---
# Column A: copy style from A82 (date column style) and set next date.
$ws.Range("A82").Copy()
$ws.Range("A83").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Cells.Item(83, 1).Value = $newRowDate

# Column BB: copy style from BB82 (default/no style) and set the value.
$ws.Range("BB82").Copy()
$ws.Range("BB83").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Cells.Item(83, 54).Value = $newForecastValue
